$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the date-serial value 45205 for every
# data row (2-353). The workbook was refreshed a day later, so bump
# every one of those cells to 45206.
$ws.Range("C2:C353").Value = 45206
